$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.60633566666667
$ws.Range("H2").Value = 52.819007
$ws.Range("I2").Value = 0.01967183396478544
$ws.Range("J2").Value = 0.01967183396478544
$ws.Range("M2").Value = 35.82261933333334
$ws.Range("N2").Value = 107.467858
$ws.Range("O2").Value = 0.1784748100644408
$ws.Range("P2").Value = 0.1784748100644408
$ws.Range("Q2").Value = 630.7050604418896
$ws.Range("R2").Value = 5676.345543977006
$ws.Range("S2").Value = 0.003510926830484298
$ws.Range("T2").Value = 0.003510926830484297
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.60633566666667
$ws.Range("H3").Value = 52.819007
$ws.Range("I3").Value = 0.01967183396478544
$ws.Range("J3").Value = 0.01967183396478544
$ws.Range("O3").Value = 0.4121780644343741
$ws.Range("P3").Value = 0.4121780644343741
$ws.Range("Q3").Value = 1456.579732165227
$ws.Range("R3").Value = 13109.21758948704
$ws.Range("S3").Value = 0.008108298447479644
$ws.Range("T3").Value = 0.00810829844747964
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.60633566666667
$ws.Range("H4").Value = 52.819007
$ws.Range("I4").Value = 0.01967183396478544
$ws.Range("J4").Value = 0.01967183396478544
$ws.Range("M4").Value = 73.74809799999998
$ws.Range("N4").Value = 221.244294
$ws.Range("O4").Value = 0.3674264480966141
$ws.Range("P4").Value = 0.3674264480966141
$ws.Range("Q4").Value = 1298.433768166228
$ws.Range("R4").Value = 11685.90391349606
$ws.Range("S4").Value = 0.007227952081227449
$ws.Range("T4").Value = 0.007227952081227447
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 17.60633566666667
$ws.Range("H5").Value = 52.819007
$ws.Range("I5").Value = 0.01967183396478544
$ws.Range("J5").Value = 0.01967183396478544
$ws.Range("M5").Value = 8.41412
$ws.Range("N5").Value = 25.24236
$ws.Range("O5").Value = 0.04192067740457094
$ws.Range("P5").Value = 0.04192067740457094
$ws.Range("Q5").Value = 148.1418210596133
$ws.Range("R5").Value = 1333.27638953652
$ws.Range("S5").Value = 0.0008246566055940523
$ws.Range("T5").Value = 0.000824656605594052
$ws.Range("I6").Value = 0.6688940299055509
$ws.Range("J6").Value = 0.6688940299055508
$ws.Range("M6").Value = 35.82261933333334
$ws.Range("N6").Value = 107.467858
$ws.Range("O6").Value = 0.1784748100644408
$ws.Range("P6").Value = 0.1784748100644408
$ws.Range("Q6").Value = 21445.62882728667
$ws.Range("R6").Value = 193010.65944558
$ws.Range("S6").Value = 0.1193807349406316
$ws.Range("T6").Value = 0.1193807349406316
$ws.Range("I7").Value = 0.6688940299055509
$ws.Range("J7").Value = 0.6688940299055508
$ws.Range("O7").Value = 0.4121780644343741
$ws.Range("P7").Value = 0.4121780644343741
$ws.Range("S7").Value = 0.2757034465581783
$ws.Range("T7").Value = 0.2757034465581783
$ws.Range("I8").Value = 0.6688940299055509
$ws.Range("J8").Value = 0.6688940299055508
$ws.Range("M8").Value = 73.74809799999998
$ws.Range("N8").Value = 221.244294
$ws.Range("O8").Value = 0.3674264480966141
$ws.Range("P8").Value = 0.3674264480966141
$ws.Range("Q8").Value = 44150.15891801887
$ws.Range("R8").Value = 397351.4302621698
$ws.Range("S8").Value = 0.2457693575612269
$ws.Range("T8").Value = 0.2457693575612269
$ws.Range("I9").Value = 0.6688940299055509
$ws.Range("J9").Value = 0.6688940299055508
$ws.Range("M9").Value = 8.41412
$ws.Range("N9").Value = 25.24236
$ws.Range("O9").Value = 0.04192067740457094
$ws.Range("P9").Value = 0.04192067740457094
$ws.Range("Q9").Value = 5037.21106346744
$ws.Range("R9").Value = 45334.89957120696
$ws.Range("S9").Value = 0.02804049084551403
$ws.Range("T9").Value = 0.02804049084551402
$ws.Range("G10").Value = 274.6625416666666
$ws.Range("H10").Value = 823.987625
$ws.Range("I10").Value = 0.3068847498029997
$ws.Range("J10").Value = 0.3068847498029996
$ws.Range("M10").Value = 35.82261933333334
$ws.Range("N10").Value = 107.467858
$ws.Range("O10").Value = 0.1784748100644408
$ws.Range("P10").Value = 0.1784748100644408
$ws.Range("Q10").Value = 9839.131675250805
$ws.Range("R10").Value = 88552.18507725725
$ws.Range("S10").Value = 0.05477119743276381
$ws.Range("T10").Value = 0.0547711974327638
$ws.Range("G11").Value = 274.6625416666666
$ws.Range("H11").Value = 823.987625
$ws.Range("I11").Value = 0.3068847498029997
$ws.Range("J11").Value = 0.3068847498029996
$ws.Range("O11").Value = 0.4121780644343741
$ws.Range("P11").Value = 0.4121780644343741
$ws.Range("Q11").Value = 22722.95036008461
$ws.Range("R11").Value = 204506.5532407615
$ws.Range("S11").Value = 0.1264911621782276
$ws.Range("T11").Value = 0.1264911621782276
$ws.Range("G12").Value = 274.6625416666666
$ws.Range("H12").Value = 823.987625
$ws.Range("I12").Value = 0.3068847498029997
$ws.Range("J12").Value = 0.3068847498029996
$ws.Range("M12").Value = 73.74809799999998
$ws.Range("N12").Value = 221.244294
$ws.Range("O12").Value = 0.3674264480966141
$ws.Range("P12").Value = 0.3674264480966141
$ws.Range("Q12").Value = 20255.84003976241
$ws.Range("R12").Value = 182302.5603578617
$ws.Range("S12").Value = 0.1127575735951343
$ws.Range("T12").Value = 0.1127575735951342
$ws.Range("G13").Value = 274.6625416666666
$ws.Range("H13").Value = 823.987625
$ws.Range("I13").Value = 0.3068847498029997
$ws.Range("J13").Value = 0.3068847498029996
$ws.Range("M13").Value = 8.41412
$ws.Range("N13").Value = 25.24236
$ws.Range("O13").Value = 0.04192067740457094
$ws.Range("P13").Value = 0.04192067740457094
$ws.Range("Q13").Value = 2311.043585088333
$ws.Range("R13").Value = 20799.392265795
$ws.Range("S13").Value = 0.01286481659687401
$ws.Range("T13").Value = 0.01286481659687401
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 4.071711
$ws.Range("H14").Value = 12.215133
$ws.Range("I14").Value = 0.004549386326664026
$ws.Range("J14").Value = 0.004549386326664025
$ws.Range("M14").Value = 35.82261933333334
$ws.Range("N14").Value = 107.467858
$ws.Range("O14").Value = 0.1784748100644408
$ws.Range("P14").Value = 0.1784748100644408
$ws.Range("Q14").Value = 145.859353188346
$ws.Range("R14").Value = 1312.734178695114
$ws.Range("S14").Value = 0.0008119508605611261
$ws.Range("T14").Value = 0.0008119508605611259
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4.071711
$ws.Range("H15").Value = 12.215133
$ws.Range("I15").Value = 0.004549386326664026
$ws.Range("J15").Value = 0.004549386326664025
$ws.Range("O15").Value = 0.4121780644343741
$ws.Range("P15").Value = 0.4121780644343741
$ws.Range("Q15").Value = 336.8544045802039
$ws.Range("R15").Value = 3031.689641221836
$ws.Range("S15").Value = 0.001875157250488585
$ws.Range("T15").Value = 0.001875157250488585
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 4.071711
$ws.Range("H16").Value = 12.215133
$ws.Range("I16").Value = 0.004549386326664026
$ws.Range("J16").Value = 0.004549386326664025
$ws.Range("M16").Value = 73.74809799999998
$ws.Range("N16").Value = 221.244294
$ws.Range("O16").Value = 0.3674264480966141
$ws.Range("P16").Value = 0.3674264480966141
$ws.Range("Q16").Value = 300.2809418556779
$ws.Range("R16").Value = 2702.528476701102
$ws.Range("S16").Value = 0.001671564859025466
$ws.Range("T16").Value = 0.001671564859025465
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 4.071711
$ws.Range("H17").Value = 12.215133
$ws.Range("I17").Value = 0.004549386326664026
$ws.Range("J17").Value = 0.004549386326664025
$ws.Range("M17").Value = 8.41412
$ws.Range("N17").Value = 25.24236
$ws.Range("O17").Value = 0.04192067740457094
$ws.Range("P17").Value = 0.04192067740457094
$ws.Range("Q17").Value = 34.25986495932
$ws.Range("R17").Value = 308.33878463388
$ws.Range("S17").Value = 0.0001907133565888486
$ws.Range("T17").Value = 0.0001907133565888486
